# Update public EPEX spot / Gaz / CO2 price workbook with the 31-jul (Prix
# Spot) and 2025-07-29 (Gaz, CO2) data points.

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column AV (31-jul) -----------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Header cell AV1 — copy AU1's formatting (bold / bordered / centered style)
# onto AV1 so it reuses the existing header style, then set its text.
$wsSpot.Range("AU1").Copy($wsSpot.Range("AV1"))
$wsSpot.Range("AV1").Value = "31-jul"

$spotValues = @(
    @{Row=2; Value=90},
    @{Row=3; Value=77.37},
    @{Row=4; Value=71.06999999999999},
    @{Row=5; Value=58.38},
    @{Row=6; Value=57.75},
    @{Row=7; Value=59.01},
    @{Row=8; Value=81.95999999999999},
    @{Row=9; Value=100.46},
    @{Row=10; Value=91.55},
    @{Row=11; Value=51.31},
    @{Row=12; Value=30.5},
    @{Row=13; Value=23.92},
    @{Row=14; Value=29.01},
    @{Row=15; Value=23.5},
    @{Row=16; Value=25.26},
    @{Row=17; Value=28.5},
    @{Row=18; Value=30},
    @{Row=19; Value=42.11},
    @{Row=20; Value=64.23999999999999},
    @{Row=21; Value=86.70999999999999},
    @{Row=22; Value=107.99},
    @{Row=23; Value=114.45},
    @{Row=24; Value=110.27},
    @{Row=25; Value=98.31999999999999}
)

foreach ($entry in $spotValues) {
    $wsSpot.Cells.Item($entry.Row, 48).Value = $entry.Value
}

# --- Sheet "Gaz": add row 45 (2025-07-29 / 33.9) ---------------------------
# Date-shaped strings get auto-recognised as real dates by Value assignment,
# so force text formatting first, then strip the formatting override back
# off (the source data stores these as plain text cells with no style).
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A45").NumberFormat = "@"
$wsGaz.Range("A45").Value = "2025-07-29"
$wsGaz.Range("A45").ClearFormats()
$wsGaz.Range("B45").Value = 33.9

# --- Sheet "CO2": add row 45 (2025-07-29 / 72.16) --------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A45").NumberFormat = "@"
$wsCo2.Range("A45").Value = "2025-07-29"
$wsCo2.Range("A45").ClearFormats()
$wsCo2.Range("B45").Value = 72.16
